# feat: add 2022-Q3 data
#
# 1) Insert a brand-new worksheet "2022-Q3" right after "总计", pushing the
#    existing quarterly sheets one slot further down the tab strip.
# 2) Populate that new sheet with the fund-holding detail rows for the
#    2022-Q3 quarter (same column layout as every other quarterly sheet).
# 3) Insert a new row 2 into "总计" (the summary sheet) carrying the
#    2022-Q3 roll-up figures, pushing the existing summary rows down and
#    renumbering the serial-index column (A) so it stays 0..7.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q3" sheet right after "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# Step 2: header row (bold, centered, thin-bordered - matches every
# other quarterly sheet's header formatting)
# ---------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q3.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------
# Step 3: fund detail rows
# columns: B=code, C=name, D=fund size, E=stock position, F=position
# ratio, G=held market value (CNY 100M), H=position rank
# ---------------------------------------------------------------------
$data = @(
  @("501092", "交银施罗德瑞思三年封闭运作混合", "51.93", "85.82", "2.07", "1.0750", 9),
  @("001764", "广发沪港深新机遇股票", "11.33", "90.34", "4.43", "0.5019", 7),
  @("009119", "广发品质回报混合A", "5.84", "92.80", "4.26", "0.2488", 8),
  @("862001", "光大阳光香港精选混合（QDII）A 人民币", "3.15", "90.62", "5.05", "0.1591", 6),
  @("862011", "光大阳光香港精选混合（QDII）A 美元", "3.15", "90.62", "5.05", "0.1591", 6),
  @("862012", "光大阳光香港精选混合（QDII）C 人民币", "3.15", "90.62", "5.05", "0.1591", 6),
  @("006671", "广发消费升级股票", "2.69", "91.34", "4.35", "0.1170", 8),
  @("005646", "中海沪港深多策略灵活配置混合", "1.20", "89.17", "5.50", "0.0660", 4),
  @("860027", "光大阳光价值30个月持有期混合B", "2.11", "91.42", "3.08", "0.0650", 10),
  @("009120", "广发品质回报混合C", "0.46", "92.80", "4.26", "0.0196", 8),
  @("860007", "光大阳光价值30个月持有期混合A", "0.59", "91.42", "3.08", "0.0182", 10)
)

$r = 2
foreach ($row in $data) {
    $q3.Cells.Item($r, 1).Value = $r - 2

    # B..G are stored as text (fund codes / percentages kept verbatim)
    $textRange = $q3.Range($q3.Cells.Item($r, 2), $q3.Cells.Item($r, 7))
    $textRange.NumberFormat = "@"
    $q3.Cells.Item($r, 2).Value = $row[0]
    $q3.Cells.Item($r, 3).Value = $row[1]
    $q3.Cells.Item($r, 4).Value = $row[2]
    $q3.Cells.Item($r, 5).Value = $row[3]
    $q3.Cells.Item($r, 6).Value = $row[4]
    $q3.Cells.Item($r, 7).Value = $row[5]

    # H is numeric (position rank)
    $q3.Cells.Item($r, 8).Value = $row[6]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 4: "总计" (summary) sheet - insert a new row 2 for 2022-Q3 and
# renumber the serial-index column for the rows pushed down
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("A2").Font.Bold = $true
$total.Range("A2").HorizontalAlignment = -4108
$total.Range("A2").VerticalAlignment = -4160
$total.Range("A2").Borders.LineStyle = 1

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 11
$total.Range("D2").Value = 2.59

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
$total.Range("A9").Value = 7
